# Add "Ancient Helmet of the Unfrozen" to the Artifacts sheet.
#
# Before:
#   row 8 "AncientHelmetOfTheUnburned" -> Base "Heavy_Stalhrim_Head", Gold 5000, Divine FALSE
#
# After:
#   row 8 "AncientHelmetOfTheUnburned" -> Base "Heavy_AncientNord_Head", Divine FALSE (Gold cleared)
#   row 9 (new) "AncientHelmetOfTheUnfrozen" -> Base "Heavy_Stalhrim_Head", Divine FALSE
#   everything that used to be row 9+ shifts down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artifacts")

# Make room for the new artifact row right above "AurielsShield" (old row 9).
$ws.Rows.Item(9).Insert()

# New row: AncientHelmetOfTheUnfrozen keeps the Stalhrim head & Divine flag that
# AncientHelmetOfTheUnburned used to have.
$ws.Cells.Item(9, 1).Value = "AncientHelmetOfTheUnfrozen"
$ws.Cells.Item(9, 2).Value = "Heavy_Stalhrim_Head"
$ws.Cells.Item(9, 6).Value = $false

# AncientHelmetOfTheUnburned (row 8) now points at the Ancient Nord head instead,
# and no longer has a Gold override.
$ws.Cells.Item(8, 2).Value = "Heavy_AncientNord_Head"
$ws.Cells.Item(8, 5).ClearContents()

# Match the author's final UI state: Artifacts tab active, selection on A10.
$ws.Activate()
$ws.Range("A10").Select()
